$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new record row at row 261 (this pushes the old rows
# 261..363 down to 262..364, matching the dimension change A1:R363 -> A1:R364).
$ws.Range("A261").EntireRow.Insert()

# Populate the newly inserted row with the new record's data. The columns
# that are constant across every data row in this sheet (A,B,C,E,F,G,H,N,Q,R)
# are copied from the surrounding rows; the record-specific columns
# (D,I,J,K,L,M,O,P) hold the new values.
$ws.Range("A261").Value = 4
$ws.Range("B261").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C261").Value = "Los Lagos"
$ws.Range("D261").Value = "07/13/2022"
$ws.Range("E261").Value = 10
$ws.Range("F261").Value = 100112023
$ws.Range("G261").Value = "Brócoli"
$ws.Range("H261").Value = "Sin especificar"
$ws.Range("I261").Value = "Primera"
$ws.Range("J261").Value = 250
$ws.Range("K261").Value = 1500
$ws.Range("L261").Value = 1500
$ws.Range("M261").Value = 1500
$ws.Range("N261").Value = "$/unidad"
$ws.Range("O261").Value = "Región Metropolitana"
$ws.Range("P261").Value = 1500
$ws.Range("Q261").Value = 1
$ws.Range("R261").Value = "Hortaliza"
